# Populate the CYM correlation matrix sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (column headers for the correlation matrix)
$ws.Range("B1").Value = "IT.CEL.SETS.P2:CYM"
$ws.Range("C1").Value = "IT.MLT.MAIN.P2:CYM"
$ws.Range("D1").Value = "SP.URB.TOTL:CYM"

# Row labels (column A)
$ws.Range("A2").Value = "IT.MLT.MAIN.P2:CYM:cor-value"
$ws.Range("A3").Value = "IT.MLT.MAIN.P2:CYM:p-value"
$ws.Range("A4").Value = "SP.POP.TOTL:CYM:cor-value"
$ws.Range("A5").Value = "SP.POP.TOTL:CYM:p-value"
$ws.Range("A6").Value = "SP.URB.TOTL:CYM:cor-value"
$ws.Range("A7").Value = "SP.URB.TOTL:CYM:p-value"

# Numeric matrix values
$ws.Range("B2").Value = -0.4908700386399798
$ws.Range("B3").Value = 0.07470030902376014

$ws.Range("B4").Value = 0.4679477595223482
$ws.Range("C4").Value = -0.9668599650002886
$ws.Range("D4").Value = 1

$ws.Range("B5").Value = 0.09152112627470861
$ws.Range("C5").Value = [double]"1.780599910272185e-08"
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = 0.4679477595223482
$ws.Range("C6").Value = -0.9668599650002886

$ws.Range("B7").Value = 0.09152112627470861
$ws.Range("C7").Value = [double]"1.780599910272185e-08"

# Styling: bold font, thin border all around, centered horizontally, top
# vertically aligned -- applied to the header/label cells. Build the style
# once on a single cell (so the style table stays compact), then copy the
# format onto the remaining label cells instead of re-issuing each
# property assignment per cell.
$base = $ws.Range("B1")
$base.Font.Bold = $true
$base.HorizontalAlignment = -4108
$base.VerticalAlignment = -4160
$base.Borders.LineStyle = 1
$base.Borders.Weight = 2

$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$ws.Range("A2:A7").PasteSpecial(-4122)
